$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "303.34"
    "E2" = "1.35%"
    "E3" = "2.95%"
    "D4" = "4.929"
    "E4" = "-3.33%"
    "D5" = "0.07835"
    "E5" = "-1.36%"
    "D6" = "2.055"
    "E6" = "-7.50%"
    "D7" = "7.838"
    "E7" = "0.86%"
    "E8" = "-0.34%"
    "D9" = "0.9211"
    "E9" = "-0.39%"
    "D10" = "0.1761"
    "E10" = "1.63%"
    "D11" = "0.07852"
    "E11" = "5.54%"
    "D12" = "0.08632"
    "E12" = "-6.81%"
    "D13" = "0.03165"
    "E13" = "4.21%"
    "E14" = "0.25%"
    "D15" = "0.001512"
    "E15" = "0.50%"
    "D16" = "0.005767"
    "E16" = "-4.50%"
    "E17" = "2,112.06%"
    "D18" = "3.469"
    "E18" = "-0.24%"
    "D19" = "2.156"
    "E19" = "-4.91%"
    "D20" = "0.3277"
    "E20" = "0.16%"
    "E21" = "0.81%"
    "D22" = "4.296"
    "E22" = "10.02%"
    "E23" = "17.22%"
    "D24" = "0.04573"
    "E24" = "-0.96%"
    "E25" = "-1.69%"
    "D26" = "0.004452"
    "E26" = "-0.45%"
    "D27" = "0.0001250"
    "E27" = "4.35%"
    "D39" = "0.01740"
    "E39" = "-1.07%"
    "D40" = "0.04793"
    "E40" = "4.15%"
    "D41" = "0.007483"
    "E41" = "7.21%"
    "D42" = "0.1363"
    "E42" = "0.23%"
    "D43" = "0.002359"
    "E43" = "7.95%"
    "D44" = "0.01060"
    "E44" = "10.78%"
    "D45" = "0.00006319"
    "E45" = "0.06%"
    "D46" = "0.00000000750"
    "E46" = "0.18%"
    "E47" = "-61.06%"
    "D48" = "0.8234"
    "E48" = "10.20%"
    "D49" = "0.00002100"
    "E49" = "0.18%"
    "D50" = "0.0002000"
    "E50" = "0.18%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"
